# Rename the single worksheet to reflect the new update date.
# Excel automatically keeps the workbook-level defined name
# ("Tandlægesystemer") and the sheet reference inside it in sync with the
# worksheet's name, so renaming here is sufficient to update both.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Opdateret d. 05-12-2025"
